# Applies the crypto-price/volume refresh described by the commit:
#   "Updated cryptos list ... with GitHub Actions"
# Only column D (Price) and column E (Volume(1h)) text values change;
# everything else in the sheet is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the workbook's inlineStr cells).
# A leading apostrophe forces Excel to keep digit-looking strings (e.g. "586.32")
# as text instead of auto-converting them to numbers; ClearFormats() then strips
# the quote-prefix cell style Excel adds so formatting stays exactly as before.
function Set-TextValue($cell, $text) {
    $cell.Value = '''' + $text
    $cell.ClearFormats()
}

$ws.Range('D2').Value = '67.665.11'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '2.488.92'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue ($ws.Range('D5')) '586.32'
$ws.Range('E5').Value = '  +0.85%  '
Set-TextValue ($ws.Range('D6')) '176.42'
$ws.Range('E6').Value = '  +4.65%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('E9').Value = '  +3.97%  '
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('E11').Value = '  +2.67%  '
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '2.949.86'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').Value = '67.529.82'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('D17').Value = '2.489.49'
$ws.Range('E17').Value = '  +0.48%  '
Set-TextValue ($ws.Range('D18')) '11.06'
$ws.Range('E18').Value = '  +0.45%  '
Set-TextValue ($ws.Range('D19')) '7.45'
$ws.Range('E19').Value = '  -0.64%  '
Set-TextValue ($ws.Range('D20')) '351.58'
$ws.Range('E20').Value = '  +0.42%  '
Set-TextValue ($ws.Range('D21')) '4.07'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('E22').Value = '  +0.17%  '
Set-TextValue ($ws.Range('D23')) '70.65'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  -0.69%  '
Set-TextValue ($ws.Range('D26')) '9.26'
$ws.Range('E26').Value = '  +1.38%  '
$ws.Range('D27').Value = '2.616.01'
$ws.Range('E27').Value = '  +0.76%  '
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '0.0₃0911'
$ws.Range('E29').Value = '  +1.34%  '
Set-TextValue ($ws.Range('D30')) '507.61'
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('E31').Value = '  +2.17%  '
$ws.Range('E32').Value = '  +2.33%  '
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('E35').Value = '  +6.56%  '
Set-TextValue ($ws.Range('D36')) '161.13'
$ws.Range('E36').Value = '  +0.92%  '
Set-TextValue ($ws.Range('D37')) '18.70'
$ws.Range('E37').Value = '  +0.22%  '
Set-TextValue ($ws.Range('D38')) '18.32'
$ws.Range('E38').Value = '  +0.37%  '
Set-TextValue ($ws.Range('D39')) '1.33'
$ws.Range('E39').Value = '  +0.89%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  +1.85%  '
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('E44').Value = '  +2.97%  '
Set-TextValue ($ws.Range('D45')) '143.49'
$ws.Range('E45').Value = '  +2.03%  '
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('E48').Value = '  +2.17%  '
$ws.Range('E49').Value = '  +0.10%  '
Set-TextValue ($ws.Range('D50')) '0.586'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('E51').Value = '  +1.78%  '
